$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Evaluation   Subjects and design" + a trailing " " run -> "Evaluation "
#    (single run), keeping the preceding "Experimental " run untouched.
# ------------------------------------------------------------------

# Locate "Evaluation   Subjects and design" (the run whose text this is).
$hit = $d.Content
$found = $hit.Find.Execute("Evaluation   Subjects and design", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "could not find 'Evaluation   Subjects and design'"
}

$keepWord = "Evaluation"
$suffixStart = $hit.Start + $keepWord.Length
$suffixEnd = $hit.End

# Mark the "   Subjects and design" tail with a throw-away direct-formatting
# flag first. This makes its formatting differ from the "Evaluation" text
# immediately to its left (and from "Experimental " further left), so the
# upcoming delete only coalesces the tail with the paragraph's trailing
# " " run -- never with "Evaluation" or "Experimental " -- leaving
# "Evaluation" untouched (its run identity / rsid survives) and producing
# a single clean "Evaluation " run once the tail (and the old separate
# trailing-space run, which shares its formatting) collapse away.
$suffix = $d.Range($suffixStart, $suffixEnd)
$suffix.Font.StrikeThrough = 1

$suffix2 = $d.Range($suffixStart, $suffixEnd)
$suffix2.Text = ""

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from just after "ABSTRACT" to just
#    after the new "Evaluation " text (Bookmarks.Add with an existing
#    name relocates it, so the old one disappears automatically).
# ------------------------------------------------------------------

# Re-find "Evaluation " now that the tail is gone, so we land right after it
# regardless of the exact character offsets.
$hit2 = $d.Content
$found2 = $hit2.Find.Execute("Evaluation ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "could not find 'Evaluation ' after edit"
}
$afterEval = $hit2.End

# Collapsed ranges are unreliable for direct bookmark placement, so drop a
# unique marker right there, locate it (non-collapsed range), anchor the
# bookmark on that, then remove the marker text again.
$anchor = $d.Range($afterEval, $afterEval)
$anchor.InsertAfter("@@GOBACKMARK@@")

$marker = $d.Content
$marker.Find.Execute("@@GOBACKMARK@@", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $marker)
$marker.Text = ""

Write-Output "done"
